# Insert a new data row right after row 109 (i.e. at row 110), pushing the
# existing rows 110-160 down to 111-161, then populate the newly inserted
# row 110 with the new weekly price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 110..160 down by one row, inheriting formatting from row 109
# (same behaviour as Excel's native "Insert Sheet Rows").
$ws.Rows.Item(110).Insert()

# Fill in the freshly inserted row 110 with the new record.
$ws.Cells.Item(110, 1).Value = 7
$ws.Cells.Item(110, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(110, 3).Value = "Ñuble"
$ws.Cells.Item(110, 4).Value = 45097
$ws.Cells.Item(110, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(110, 5).Value = 16
$ws.Cells.Item(110, 6).Value = 100112037
$ws.Cells.Item(110, 7).Value = "Cebollín"
$ws.Cells.Item(110, 8).Value = "Sin especificar"
$ws.Cells.Item(110, 9).Value = "Primera"
$ws.Cells.Item(110, 10).Value = 200
$ws.Cells.Item(110, 11).Value = 6000
$ws.Cells.Item(110, 12).Value = 7000
$ws.Cells.Item(110, 13).Value = 6500
$ws.Cells.Item(110, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(110, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(110, 16).Value = 181
$ws.Cells.Item(110, 17).Value = 36
$ws.Cells.Item(110, 18).Value = "Hortaliza"
